$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns that actually carry (non-empty) data in rows 2-5 of this sheet.
# NOTE: a few columns are intentionally excluded from the swap because their
# value is identical across all four rows, so swapping is a no-op, and
# touching them only risks incidental side effects on write-back:
#  - "Y" (Startdatum) / "AA" (Slutdatum): identical date-text "2023-09-10"
#    in every row; writing it back risks Excel auto-converting the
#    date-looking text into a date serial number.
#  - "I", "K", "AT", "AY": empty in every row; writing an empty string back
#    causes the cell to be dropped from the saved file.
$cols = @("A","B","C","D","E","F","G","H","P","Q","R","S","T","U","V","W","Z","AB","AD","AE","AG","AW","AX")

function Swap-Rows($r1, $r2) {
    foreach ($col in $cols) {
        $cell1 = $ws.Range($col + $r1)
        $cell2 = $ws.Range($col + $r2)
        $v1 = $cell1.Value2
        $v2 = $cell2.Value2

        $cell1.Value2 = $v2
        $cell2.Value2 = $v1
    }
}

Swap-Rows 2 4
Swap-Rows 3 5
